# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" (E) and "Valor Mora" (F) values between row 16 and row 18.
# Row 16 previously: 2310 / 32480  -> becomes 2312 / 46400
# Row 18 previously: 2312 / 46400  -> becomes 2310 / 32480
# Row 17 (2311 / 46400) stays unchanged.
$ws.Range("E16").Value = "2312"
$ws.Range("F16").Value = 46400

$ws.Range("E18").Value = "2310"
$ws.Range("F18").Value = 32480
